$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- New shared strings used below ----
$sProcesoId   = "Proceso ID"
$sTipoTrans   = "Tipo de Transición"
$sTiempoMs    = "Tiempo (ms)"
$sReadyRun    = "Ready → Running"
$sTiempoProm  = "Tiempo promedio"
$sNota9       = "Esto fue de simulacion para la Medir tiempos de transición entre estado"

# ---- Apply box-outline borders (left/right/top/bottom on the outer edge only) ----
function Set-BoxBorder($range) {
    $range.Borders.Item(7).LineStyle = 1   # xlEdgeLeft
    $range.Borders.Item(7).Weight = 2      # xlThin
    $range.Borders.Item(10).LineStyle = 1  # xlEdgeRight
    $range.Borders.Item(10).Weight = 2
    $range.Borders.Item(8).LineStyle = 1   # xlEdgeTop
    $range.Borders.Item(8).Weight = 2
    $range.Borders.Item(9).LineStyle = 1   # xlEdgeBottom
    $range.Borders.Item(9).Weight = 2
}

# ---- Apply a full per-cell border (left+right+top+bottom on every cell) ----
function Set-FullGridBorder($range) {
    $range.Borders.Item(7).LineStyle = 1
    $range.Borders.Item(7).Weight = 2
    $range.Borders.Item(10).LineStyle = 1
    $range.Borders.Item(10).Weight = 2
    $range.Borders.Item(8).LineStyle = 1
    $range.Borders.Item(8).Weight = 2
    $range.Borders.Item(9).LineStyle = 1
    $range.Borders.Item(9).Weight = 2
    $range.Borders.Item(11).LineStyle = 1  # xlInsideVertical
    $range.Borders.Item(11).Weight = 2
    $range.Borders.Item(12).LineStyle = 1  # xlInsideHorizontal
    $range.Borders.Item(12).Weight = 2
}

# Row 8 gets a box border first (matches original authoring order)
Set-BoxBorder($ws.Range("A8:D8"))

# Rows 2:5 get a box border around the whole block
Set-BoxBorder($ws.Range("A2:D5"))

# Row 6 (existing merged title) gets a box border too
Set-BoxBorder($ws.Range("A6:D6"))

# ---- New row 9: merged subtitle note ----
$ws.Range("A9:D9").Merge()
$ws.Range("A9").Value = $sNota9
$a9 = $ws.Range("A9")
$a9.Font.Bold = $true
$a9.Font.Italic = $true
$a9.Font.Underline = 2   # xlUnderlineStyleSingle
$a9.HorizontalAlignment = -4108  # xlCenter
Set-BoxBorder($ws.Range("A9:D9"))

# ---- New row 10: headers for the second table ----
$ws.Range("A10").Value = $sProcesoId
$ws.Range("B10").Value = $sTipoTrans
$ws.Range("C10").Value = $sTiempoMs
Set-FullGridBorder($ws.Range("A10:C10"))

# ---- New rows 11-15: transition data ----
$ws.Range("A11").Value = 1
$ws.Range("B11").Value = $sReadyRun
$ws.Range("C11").Value = 1661

$ws.Range("A12").Value = 2
$ws.Range("B12").Value = $sReadyRun
$ws.Range("C12").Value = 4720

$ws.Range("A13").Value = 3
$ws.Range("B13").Value = $sReadyRun
$ws.Range("C13").Value = 3023

$ws.Range("A14").Value = 4
$ws.Range("B14").Value = $sReadyRun
$ws.Range("C14").Value = 2038

$ws.Range("A15").Value = 5
$ws.Range("B15").Value = $sReadyRun
$ws.Range("C15").Value = 1680

# ---- New row 16: average ----
$ws.Range("B16").Value = $sTiempoProm
$ws.Range("C16").Value = 2624

Set-BoxBorder($ws.Range("A11:C16"))

# Number format (thousands separator, 0 decimals) for the Tiempo (ms) column
$ws.Range("C11:C16").NumberFormat = "#,##0"

# ---- Selection & view ----
$ws.Range("E10").Select()
